$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; B=86; C="house/house019.jpg"; D="küssen"; E="house" },
    @{ Row=3; B=29; C="house/house022.jpg"; D="süßen"; E="house" },
    @{ Row=4; B=106; C="dog/dog030.jpg"; D="wachsen"; E="dog" },
    @{ Row=5; B=54; C="house/house007.jpg"; D="öffnen"; E="house" },
    @{ Row=6; B=45; C="house/house008.jpg"; D="holen"; E="house" },
    @{ Row=7; B=120; C="house/house015.jpg"; D="binden"; E="house" },
    @{ Row=8; B=67; C="house/house021.jpg"; D="atmen"; E="house" },
    @{ Row=9; B=84; C="house/house011.jpg"; D="narren"; E="house" },
    @{ Row=10; B=56; C="dog/dog006.jpg"; D="prüfen"; E="dog" },
    @{ Row=11; B=71; C="dog/dog007.jpg"; D="legen"; E="dog" },
    @{ Row=12; B=91; C="dog/dog001.jpg"; D="betteln"; E="dog" },
    @{ Row=13; B=100; C="dog/dog028.jpg"; D="sparen"; E="dog" },
    @{ Row=14; B=22; C="house/house023.jpg"; D="deuten"; E="house" },
    @{ Row=15; B=28; C="dog/dog004.jpg"; D="tollen"; E="dog" },
    @{ Row=16; B=110; C="house/house024.jpg"; D="quellen"; E="house" },
    @{ Row=17; B=108; C="dog/dog023.jpg"; D="lächeln"; E="dog" },
    @{ Row=18; B=16; C="house/house004.jpg"; D="stoppen"; E="house" },
    @{ Row=19; B=89; C="dog/dog002.jpg"; D="meinen"; E="dog" },
    @{ Row=20; B=48; C="dog/dog031.jpg"; D="achten"; E="dog" },
    @{ Row=21; B=116; C="house/house001.jpg"; D="duschen"; E="house" },
    @{ Row=22; B=61; C="house/house026.jpg"; D="rechnen"; E="house" },
    @{ Row=23; B=102; C="dog/dog020.jpg"; D="ärgern"; E="dog" },
    @{ Row=24; B=93; C="house/house002.jpg"; D="zögern"; E="house" },
    @{ Row=25; B=27; C="dog/dog011.jpg"; D="herrschen"; E="dog" },
    @{ Row=26; B=119; C="dog/dog018.jpg"; D="danken"; E="dog" },
    @{ Row=27; B=59; C="house/house029.jpg"; D="reisen"; E="house" },
    @{ Row=28; B=4; C="house/house005.jpg"; D="heißen"; E="house" },
    @{ Row=29; B=111; C="dog/dog003.jpg"; D="piepen"; E="dog" },
    @{ Row=30; B=26; C="dog/dog025.jpg"; D="hassen"; E="dog" },
    @{ Row=31; B=30; C="dog/dog008.jpg"; D="wecken"; E="dog" },
    @{ Row=32; B=103; C="house/house006.jpg"; D="lassen"; E="house" },
    @{ Row=33; B=117; C="dog/dog009.jpg"; D="mögen"; E="dog" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
